# Add 135 new rows (A430:D564) for 'Exportação' ('Exportação' is a new
# category under Variavel, column C) across the 9 countries that have
# export data (all of the top-10 countries already in the sheet except
# Haiti), years 2008-2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$exportData = @(
    @(2008, 'China', 5400000),
    @(2009, 'China', 1500000),
    @(2010, 'China', 1500000),
    @(2011, 'China', 1900000),
    @(2012, 'China', 2000000),
    @(2013, 'China', 1900000),
    @(2014, 'China', 3700000),
    @(2015, 'China', 8200000),
    @(2016, 'China', 10000000),
    @(2017, 'China', 9400000),
    @(2018, 'China', 6300000),
    @(2019, 'China', 3600000),
    @(2020, 'China', 1600000),
    @(2021, 'China', 4200000),
    @(2022, 'China', 2900000),
    @(2008, 'Alemanha', 358000000),
    @(2009, 'Alemanha', 355700000),
    @(2010, 'Alemanha', 392900000),
    @(2011, 'Alemanha', 414500000),
    @(2012, 'Alemanha', 397300000),
    @(2013, 'Alemanha', 390300000),
    @(2014, 'Alemanha', 367500000),
    @(2015, 'Alemanha', 392800000),
    @(2016, 'Alemanha', 360900000),
    @(2017, 'Alemanha', 383600000),
    @(2018, 'Alemanha', 373100000),
    @(2019, 'Alemanha', 383900000),
    @(2020, 'Alemanha', 366200000),
    @(2021, 'Alemanha', 368900000),
    @(2022, 'Alemanha', 353300000),
    @(2008, 'Japão', 400000),
    @(2009, 'Japão', 400000),
    @(2010, 'Japão', 200000),
    @(2011, 'Japão', 300000),
    @(2012, 'Japão', 200000),
    @(2013, 'Japão', 400000),
    @(2014, 'Japão', 200000),
    @(2015, 'Japão', 300000),
    @(2016, 'Japão', 200000),
    @(2017, 'Japão', 200000),
    @(2018, 'Japão', 200000),
    @(2019, 'Japão', 100000),
    @(2020, 'Japão', 200000),
    @(2021, 'Japão', 300000),
    @(2022, 'Japão', 300000),
    @(2008, 'Paises Baixos', 17700000),
    @(2009, 'Paises Baixos', 16700000),
    @(2010, 'Paises Baixos', 14900000),
    @(2011, 'Paises Baixos', 24400000),
    @(2012, 'Paises Baixos', 28200000),
    @(2013, 'Paises Baixos', 25800000),
    @(2014, 'Paises Baixos', 31300000),
    @(2015, 'Paises Baixos', 37600000),
    @(2016, 'Paises Baixos', 66800000),
    @(2017, 'Paises Baixos', 79500000),
    @(2018, 'Paises Baixos', 72700000),
    @(2019, 'Paises Baixos', 64900000),
    @(2020, 'Paises Baixos', 96000000),
    @(2021, 'Paises Baixos', 116500000),
    @(2022, 'Paises Baixos', 112300000),
    @(2008, 'Paraguai', 0),
    @(2009, 'Paraguai', 0),
    @(2010, 'Paraguai', 0),
    @(2011, 'Paraguai', 0),
    @(2012, 'Paraguai', 0),
    @(2013, 'Paraguai', 0),
    @(2014, 'Paraguai', 0),
    @(2015, 'Paraguai', 0),
    @(2016, 'Paraguai', 0),
    @(2017, 'Paraguai', 0),
    @(2018, 'Paraguai', 100000),
    @(2019, 'Paraguai', 0),
    @(2020, 'Paraguai', 0),
    @(2021, 'Paraguai', 0),
    @(2022, 'Paraguai', 0),
    @(2008, 'Espanha', 1691400000),
    @(2009, 'Espanha', 1460700000),
    @(2010, 'Espanha', 1715600000),
    @(2011, 'Espanha', 2203100000),
    @(2012, 'Espanha', 2141100000),
    @(2013, 'Espanha', 1845000000),
    @(2014, 'Espanha', 2308800000),
    @(2015, 'Espanha', 2439500000),
    @(2016, 'Espanha', 2262600000),
    @(2017, 'Espanha', 2331100000),
    @(2018, 'Espanha', 2030800000),
    @(2019, 'Espanha', 2138700000),
    @(2020, 'Espanha', 2017200000),
    @(2021, 'Espanha', 2298500000),
    @(2022, 'Espanha', 2089500000),
    @(2008, 'Reino Unido', 43300000),
    @(2009, 'Reino Unido', 54500000),
    @(2010, 'Reino Unido', 89800000),
    @(2011, 'Reino Unido', 86600000),
    @(2012, 'Reino Unido', 80000000),
    @(2013, 'Reino Unido', 95100000),
    @(2014, 'Reino Unido', 103700000),
    @(2015, 'Reino Unido', 94200000),
    @(2016, 'Reino Unido', 78600000),
    @(2017, 'Reino Unido', 96200000),
    @(2018, 'Reino Unido', 114400000),
    @(2019, 'Reino Unido', 95400000),
    @(2020, 'Reino Unido', 88300000),
    @(2021, 'Reino Unido', 35800000),
    @(2022, 'Reino Unido', 32600000),
    @(2008, 'Estados Unidos', 463800000),
    @(2009, 'Estados Unidos', 397400000),
    @(2010, 'Estados Unidos', 400900000),
    @(2011, 'Estados Unidos', 416500000),
    @(2012, 'Estados Unidos', 400400000),
    @(2013, 'Estados Unidos', 414700000),
    @(2014, 'Estados Unidos', 404700000),
    @(2015, 'Estados Unidos', 418900000),
    @(2016, 'Estados Unidos', 379000000),
    @(2017, 'Estados Unidos', 345400000),
    @(2018, 'Estados Unidos', 348700000),
    @(2019, 'Estados Unidos', 357300000),
    @(2020, 'Estados Unidos', 360400000),
    @(2021, 'Estados Unidos', 326600000),
    @(2022, 'Estados Unidos', 278500000),
    @(2008, 'Uruguai', 13400000),
    @(2009, 'Uruguai', 2000000),
    @(2010, 'Uruguai', 2600000),
    @(2011, 'Uruguai', 2400000),
    @(2012, 'Uruguai', 19700000),
    @(2013, 'Uruguai', 15200000),
    @(2014, 'Uruguai', 4300000),
    @(2015, 'Uruguai', 5800000),
    @(2016, 'Uruguai', 3400000),
    @(2017, 'Uruguai', 5800000),
    @(2018, 'Uruguai', 18500000),
    @(2019, 'Uruguai', 5800000),
    @(2020, 'Uruguai', 6900000),
    @(2021, 'Uruguai', 5300000),
    @(2022, 'Uruguai', 5300000)
)

$startRow = 430
for ($i = 0; $i -lt $exportData.Count; $i++) {
    $row = $startRow + $i
    $entry = $exportData[$i]
    $year = $entry[0]
    $country = $entry[1]
    $liters = $entry[2]

    $ws.Cells.Item($row, 1).Value = $year
    $ws.Cells.Item($row, 2).Value = $country
    $ws.Cells.Item($row, 3).Value = "Exportação"

    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $liters
    $cell.NumberFormat = "#,##0.00"
}

$lastRow = $startRow + $exportData.Count - 1

# Select the newly added block, matching the post-edit sheetView selection.
$ws.Range("A" + $startRow + ":D" + $lastRow).Select()

# Best-effort: scroll the viewport so the new rows are visible (mirrors the
# `topLeftCell="A456"` attribute Excel recorded after the paste). Some hosts
# may not persist this into the saved sheetView without frozen panes.
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 456
    $win.ScrollColumn = 1
}
